# Updates cryptos list: refreshes price (col D) and volume/1h change (col E)
# values, and swaps the row positions of three coin pairs (Chainlink /
# WrappedBTC at rows 17-18, and Stacks / EnergySwap at rows 50-51),
# matching commit: 'Updated cryptos list on Sat Jun  1 09:39:19 UTC 2024
# with GitHub Actions'.
#
# Column D cells get NumberFormat forced to text ("@") before the value
# is set, because some of the new price strings (e.g. '595.59', '1.00')
# look numeric and Excel's COM layer would otherwise silently convert
# them to real numbers, whereas the source data stores them as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.686.37'
$ws.Range("E2").Value = '  -0.35%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.791.40'
$ws.Range("E3").Value = '  +1.43%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '595.59'
$ws.Range("E5").Value = '  +0.61%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.20'
$ws.Range("E6").Value = '  +0.30%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.790.56'
$ws.Range("E7").Value = '  +1.37%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.520'
$ws.Range("E9").Value = '  +0.03%  '

$ws.Range("E10").Value = '  +0.38%  '

$ws.Range("E11").Value = '  -1.54%  '

$ws.Range("E12").Value = '  +0.10%  '

$ws.Range("E13").Value = '  -1.23%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.11'
$ws.Range("E14").Value = '  +0.51%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.426.67'
$ws.Range("E15").Value = '  +1.47%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.783.82'
$ws.Range("E16").Value = '  +1.92%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.629.27'
$ws.Range("E17").Value = '  -0.38%  '

$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.40'
$ws.Range("E18").Value = '  +2.99%  '

$ws.Range("E19").Value = '  +0.64%  '

$ws.Range("E20").Value = '  +0.07%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.99'
$ws.Range("E21").Value = '  -5.99%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '458.67'
$ws.Range("E22").Value = '  -1.10%  '

$ws.Range("E23").Value = '  +0.38%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000155'
$ws.Range("E24").Value = '  +6.14%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.28'
$ws.Range("E25").Value = '  -0.44%  '

$ws.Range("E26").Value = '  +1.97%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.12'

$ws.Range("E28").Value = '  -0.01%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.01'
$ws.Range("E29").Value = '  -0.10%  '

$ws.Range("E30").Value = '  +0.18%  '

$ws.Range("E31").Value = '  +4.33%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.24'
$ws.Range("E32").Value = '  +0.31%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.63'
$ws.Range("E33").Value = '  -0.13%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  -0.11%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.08'
$ws.Range("E35").Value = '  -0.03%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.728.87'
$ws.Range("E36").Value = '  +1.02%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0999'
$ws.Range("E37").Value = '  -0.44%  '

$ws.Range("E38").Value = '  -1.61%  '

$ws.Range("E39").Value = '  +0.27%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.993'
$ws.Range("E40").Value = '  -0.07%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.77'
$ws.Range("E41").Value = '  +0.43%  '

$ws.Range("E42").Value = '  -0.08%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '44.21'
$ws.Range("E44").Value = '  +0.43%  '

$ws.Range("E45").Value = '  +2.73%  '

$ws.Range("E46").Value = '  +0.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '149.43'
$ws.Range("E47").Value = '  +3.38%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.28'
$ws.Range("E48").Value = '  -1.25%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '392.76'
$ws.Range("E49").Value = '  +1.52%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '26.59'
$ws.Range("E50").Value = '  +6.85%  '

$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.82'
$ws.Range("E51").Value = '  -4.10%  '

